$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(100, 101, 102, 108, 109, 116, 117, 118, 119, 126, 127, 128, 134, 137, 138, 139, 141)

foreach ($r in $rows) {
    $ws.Range("B$r").Value = "train"
}
